# Generate Report for Handoff
# Updates the localization-status workbook so that the previous handoff
# file GUID (e960e1ad-b4d9-4e3d-86cf-82409db988d2) is replaced by the new
# one (adc849d1-f7db-4533-a844-66c4cb929d9d), along with the associated
# content hash and timestamps, across the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$oldGuid = "e960e1ad-b4d9-4e3d-86cf-82409db988d2"
$newGuid = "adc849d1-f7db-4533-a844-66c4cb929d9d"
$oldHash = "f388c5d9ec894f07404001512cf78f153a920329"
$newHash = "fc8db7f6ab9729159ffa3c4cae9c7047352a7372"

$ghBase = "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/2b084e0477e002a573eec66f946769c9010d0cc0/e2e/"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2017-02-21 04:46:01"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "$ghBase$oldGuid.md", "", "", "e2e\$newGuid.md")

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2017-02-21 04:45:45"

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "$ghBase$oldGuid.md", "", "", "$newGuid.md")

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2017-02-21 04:46:01"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "$ghBase$oldGuid.md", "", "", "$newGuid.md")
